# Update the responsible-team text for SPRINT 1 (row 2) and SPRINT 3 (row 9)
# from "Gabriel, Lucas, Gabriela, Higor" to "Gabriel, Gabriela, Higor".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("D2").Value = "Gabriel, Gabriela, Higor"
$ws.Range("D9").Value = "Gabriel, Gabriela, Higor"

# Move the active selection to F16, matching the saved view state.
$ws.Range("F16").Select()
